# Fix: dates in D2/E2 were stored as text (shared strings "2/15/2025",
# "3/20/2025") instead of real date values, and the sheet's saved
# selection pointed at H5 (outside the actual data), which caused an
# "infinite view" in the listing page. This applies real Excel date
# serial values with a date number format, and resets the selection to
# the top of the data (D2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# openingDate (D2) -> 2/15/2025, closingDate (E2) -> 3/20/2025
$ws.Range("D2").Value = 45703
$ws.Range("E2").Value = 45736

# Apply a real date number format (built-in date format id 14) to D2,
# then copy that exact formatting onto E2 so both cells share one style.
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reset the lingering selection (was H5) back onto the data.
$ws.Range("D2").Select() | Out-Null
